$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.097.54'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').Value = '1.808.81'
$ws.Range('E3').Value = '  -2.15%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').Value = '''232.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = '''0.612'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').Value = '''40.54'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.06%  '
$ws.Range('D9').Value = '''0.325'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.47%  '
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('D11').Value = '''0.0999'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').Value = '2.071.35'
$ws.Range('E12').Value = '  -2.11%  '
$ws.Range('D13').Value = '1.808.82'
$ws.Range('E13').Value = '  -2.18%  '
$ws.Range('D14').Value = '''0.664'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').Value = '''11.07'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.95%  '
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('D17').Value = '35.060.40'
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('D18').Value = '''69.66'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('D19').Value = '0.0₃0789'
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D20').Value = '''237.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.36%  '
$ws.Range('E21').Value = '  -2.21%  '
$ws.Range('E22').Value = '  -1.04%  '
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').Value = '''2.24'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.51%  '
$ws.Range('D25').Value = '''171.78'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').Value = '''7.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('D27').Value = '''17.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.09%  '
$ws.Range('D29').Value = '''1.57'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +17.90%  '
$ws.Range('D30').Value = '''1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('E31').Value = '  +4.18%  '
$ws.Range('E32').Value = '  +3.94%  '
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('D34').Value = '''1.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.07%  '
$ws.Range('D35').Value = '''0.698'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.03%  '
$ws.Range('E36').Value = '  +5.64%  '
$ws.Range('D37').Value = '''92.57'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.05%  '
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').Value = '1.315.26'
$ws.Range('E39').Value = '  -1.82%  '
$ws.Range('E40').Value = '  -0.33%  '
$ws.Range('D41').Value = '''0.995'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.24%  '
$ws.Range('E42').Value = '  +0.69%  '
$ws.Range('D43').Value = '''14.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.22%  '
$ws.Range('E44').Value = '  -6.83%  '
$ws.Range('E45').Value = '  -2.26%  '
$ws.Range('D46').Value = '''6.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.02%  '
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('D48').Value = '1.988.76'
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('D49').Value = '''1.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('D50').Value = '''0.0670'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.35%  '
$ws.Range('D51').Value = '''99.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.59%  '
